# The workbook originally has a stray first row ("ID", "Entidade", "Data
# Venc.", "Categoria", "Dias", "Valor Pendente") sitting above the real
# header row (blank, "Entidade", "Categoria", "Dias", "Valor Pendente").
# Remove that stray first row so the real header becomes row 1, widen
# column E to fit the now-wider "Valor Pendente" values, and reset the
# selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

$ws.Columns("E").ColumnWidth = 15

$ws.Range("A1").Select()
